# Insert one new data row at sheet row 717 (pushing the existing
# 2026/12/29 .. 2027/01/05 rows down by one), and fill the new row with
# the 2026/01/26 08:00-ish entry that was missing.
#
# Equivalent data-level effect: the row that used to be at position 717
# ("2026/12/29", "火", 13, 201) now lives at 718, and so on through the
# end of the sheet (758 -> 759); a brand-new row 717 holds
# ("2026/01/26", "月", 13, 201). dimension grows from A1:D758 to A1:D759.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 717 and everything below it down by one row.
$ws.Rows("717:717").Insert()

# Populate the newly-inserted (now blank) row 717.
# Force column A to stay plain text ("2026/01/26") instead of being
# auto-parsed into a date serial, matching how the rest of the date
# column is stored (inline literal text, default "General" style).
$ws.Range("A717").NumberFormat = "@"
$ws.Range("A717").Value = "2026/01/26"
$ws.Range("A717").Style = "Normal"

$ws.Range("B717").Value = "月"
$ws.Range("C717").Value = 13
$ws.Range("D717").Value = 201
